$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Clientes
$ws2 = $wb.Worksheets.Item(2)   # Dados Manjerico

# --- Sheet "Dados Manjerico": Manjerico's own street address no longer
# duplicates the house number (it already lives in its own "Numero" column) ---
$ws2.Range("C2").Value = "Rua Itacuruçá"

# --- Add a new "E-mail" column (K) for Manjerico, header styled like the
# other sheet's special bold-black header, value as a mailto hyperlink ---
$ws2.Range("K1").Value = "E-mail"
$ws2.Range("K1").Font.Bold = $true
$ws2.Range("K1").Font.Color = 0

$ws2.Range("K2").Value = "manjerico@manjerico.com.br"
$ws2.Hyperlinks.Add($ws2.Range("K2"), "mailto:manjerico@manjerico.com.br") | Out-Null

# --- Make "Dados Manjerico" (now holding the sender's data) the active /
# selected sheet, replacing "Clientes" ---
$ws2.Range("K1:K2").Select()
$ws2.Activate()

$wb.Save()
